# feat: add 2022-Q1 data
#
# The workbook's physical sheet that used to hold the "总计" (grand total)
# summary table is renamed to "2022-Q1" and repurposed to hold the new
# quarter's per-fund breakdown (same shape as the other quarterly sheets).
# A brand-new "总计" sheet is appended at the end, holding the old summary
# table with one new leading row ("2022-Q1") and the index column
# renumbered.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Style/reference helpers.
#   - $hdrRef   : a header cell (s=2 -> bold+border, centered) to clone
#                 formatting from.
#   - $idxRef   : an index-column cell (s=2, same style index actually)
#                 to clone formatting from.
#   - $plainRef : a cell with NO explicit style (style 0 / default).
# All three come from the still-untouched "2021-Q4" sheet so they are not
# disturbed by the edits made below.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$hdrRef   = $refSheet.Range("B1")
$plainRef = $refSheet.Range("B2")

# Writes $val into $rng as literal TEXT (never auto-coerced to a number),
# while keeping the cell's style at "no explicit style" (matches the
# plain data cells in the source sheets, which carry no `s` attribute).
function Set-PlainText($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $plainRef.Copy()
    $rng.PasteSpecial(-4122)   # xlPasteFormats: formatting only, value untouched
}

# Writes $val into $rng as literal TEXT, using the header style (s=2).
function Set-HeaderText($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $hdrRef.Copy()
    $rng.PasteSpecial(-4122)
}

# Writes a plain number into $rng and applies the header/index style (s=2).
function Set-StyledNumber($rng, $val) {
    $rng.Value = $val
    $hdrRef.Copy()
    $rng.PasteSpecial(-4122)
}

# =======================================================================
# 1) Rename "总计" -> "2022-Q1" and replace its content with the new
#    per-fund breakdown for 2022-Q1.
# =======================================================================
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

Set-HeaderText $q1.Range("B1") "基金代码"
Set-HeaderText $q1.Range("C1") "基金名称"
Set-HeaderText $q1.Range("D1") "基金规模"
Set-HeaderText $q1.Range("E1") "股票总仓位"
Set-HeaderText $q1.Range("F1") "仓位占比"
Set-HeaderText $q1.Range("G1") "持有市值(亿元)"
Set-HeaderText $q1.Range("H1") "仓位排名"

$q1Rows = @(
    @("539003", "建信富时100指数（QDII）人民币A", "0.71", "92.86", "4.30", "0.0305", 9),
    @("008707", "建信富时100指数（QDII）美元现汇A", "0.71", "92.86", "4.30", "0.0305", 9),
    @("008706", "建信富时100指数（QDII）人民币C", "0.20", "92.86", "4.30", "0.0086", 9),
    @("008708", "建信富时100指数（QDII）美元现汇C", "0.20", "92.86", "4.30", "0.0086", 9)
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]
    Set-StyledNumber $q1.Range("A$r") $i
    Set-PlainText $q1.Range("B$r") $row[0]
    Set-PlainText $q1.Range("C$r") $row[1]
    Set-PlainText $q1.Range("D$r") $row[2]
    Set-PlainText $q1.Range("E$r") $row[3]
    Set-PlainText $q1.Range("F$r") $row[4]
    Set-PlainText $q1.Range("G$r") $row[5]
    $q1.Range("H$r").Value = $row[6]
}

$excel.CutCopyMode = $false

# =======================================================================
# 2) Append a brand-new "总计" sheet (after "2022-Q1") holding the grand
#    total table: the previous rows plus a new leading "2022-Q1" row,
#    with the index column renumbered 0..5.
# =======================================================================
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totRef = $wb.Worksheets.Item("2021-Q4")
$hdrRef2 = $totRef.Range("B1")

function Set-HeaderText2($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $hdrRef2.Copy()
    $rng.PasteSpecial(-4122)
}
function Set-StyledNumber2($rng, $val) {
    $rng.Value = $val
    $hdrRef2.Copy()
    $rng.PasteSpecial(-4122)
}

Set-HeaderText2 $total.Range("B1") "日期"
Set-HeaderText2 $total.Range("C1") "持有数量(只)"
Set-HeaderText2 $total.Range("D1") "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 4, 0.08),
    @("2021-Q4", 4, 0.23),
    @("2021-Q3", 4, 0.07000000000000001),
    @("2021-Q2", 4, 0.11),
    @("2021-Q1", 4, 0.1),
    @("2020-Q4", 4, 0.11)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    Set-StyledNumber2 $total.Range("A$r") $i
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
}

$excel.CutCopyMode = $false
